$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 48:50 (remove last 3 data rows)
$ws.Rows("48:50").Delete()

# Update B2:C47 values
$ws.Range("B2").Value = 1.182783223950844
$ws.Range("C2").Value = 0.588493076750093
$ws.Range("B3").Value = 6.267955951194165
$ws.Range("C3").Value = 1.439892760547884
$ws.Range("B4").Value = 17.99679854773782
$ws.Range("C4").Value = 2.009086801263544
$ws.Range("B5").Value = 24.6241679219543
$ws.Range("C5").Value = 3.209923369010773
$ws.Range("B6").Value = 25.64548164713591
$ws.Range("C6").Value = 3.967170290534838
$ws.Range("B7").Value = 28.78439610034238
$ws.Range("C7").Value = 4.818625260531809
$ws.Range("B8").Value = 29.06193199555131
$ws.Range("C8").Value = 5.615194114720508
$ws.Range("B9").Value = 29.59465251538218
$ws.Range("C9").Value = 6.344054604464469
$ws.Range("B10").Value = 29.90179677752234
$ws.Range("C10").Value = 7.491534512736186
$ws.Range("B11").Value = 31.2870511541718
$ws.Range("C11").Value = 8.155531395671357
$ws.Range("B12").Value = 32.37159994294669
$ws.Range("C12").Value = 9.242282652153332
$ws.Range("B13").Value = 32.46151952625841
$ws.Range("C13").Value = 9.909265447822202
$ws.Range("B14").Value = 32.54294375292982
$ws.Range("C14").Value = 10.33839816564545
$ws.Range("B15").Value = 35.20062906453335
$ws.Range("C15").Value = 11.63782164517427
$ws.Range("B16").Value = 36.64645947817002
$ws.Range("C16").Value = 12.30629773594152
$ws.Range("B17").Value = 39.33688517122047
$ws.Range("C17").Value = 12.86714603728094
$ws.Range("B18").Value = 40.20959427425738
$ws.Range("C18").Value = 13.54475126136632
$ws.Range("B19").Value = 43.58924743161049
$ws.Range("C19").Value = 14.1915518659535
$ws.Range("B20").Value = 45.3776442754611
$ws.Range("C20").Value = 15.1705662208524
$ws.Range("B21").Value = 45.57123167624007
$ws.Range("C21").Value = 15.84872914879021
$ws.Range("B22").Value = 47.72028453017151
$ws.Range("C22").Value = 16.436848907948
$ws.Range("B23").Value = 48.07404658882837
$ws.Range("C23").Value = 17.44086136674384
$ws.Range("B24").Value = 53.24798226570667
$ws.Range("C24").Value = 18.04963291163811
$ws.Range("B25").Value = 54.10847605624152
$ws.Range("C25").Value = 19.30505671579645
$ws.Range("B26").Value = 57.12810609209131
$ws.Range("C26").Value = 20.11047165871832
$ws.Range("B27").Value = 57.28873318221724
$ws.Range("C27").Value = 20.96713294172871
$ws.Range("B28").Value = 58.27719513929562
$ws.Range("C28").Value = 21.81332224797398
$ws.Range("B29").Value = 62.31644247555465
$ws.Range("C29").Value = 22.37289868520194
$ws.Range("B30").Value = 62.44304949078827
$ws.Range("C30").Value = 22.94119909625795
$ws.Range("B31").Value = 64.88795701797191
$ws.Range("C31").Value = 23.59446120042404
$ws.Range("B32").Value = 65.08437267076951
$ws.Range("C32").Value = 24.59578597423474
$ws.Range("B33").Value = 68.86586574388853
$ws.Range("C33").Value = 25.11423832814771
$ws.Range("B34").Value = 69.70216695675539
$ws.Range("C34").Value = 26.14155177784091
$ws.Range("B35").Value = 69.79535098445156
$ws.Range("C35").Value = 26.87292741403455
$ws.Range("B36").Value = 70.45759044980939
$ws.Range("C36").Value = 27.63973466582038
$ws.Range("B37").Value = 70.59504811641017
$ws.Range("C37").Value = 28.67482152301453
$ws.Range("B38").Value = 73.7493645851842
$ws.Range("C38").Value = 29.57772640018024
$ws.Range("B39").Value = 74.84085334728266
$ws.Range("C39").Value = 30.43066485844317
$ws.Range("B40").Value = 77.29061552584668
$ws.Range("C40").Value = 31.04249034196802
$ws.Range("B41").Value = 77.53481435264875
$ws.Range("C41").Value = 31.96554655960691
$ws.Range("B42").Value = 86.45128086456569
$ws.Range("C42").Value = 32.54456371191696
$ws.Range("B43").Value = 86.51931084564436
$ws.Range("C43").Value = 33.2505875237318
$ws.Range("B44").Value = 86.69619166932952
$ws.Range("C44").Value = 33.93648194340026
$ws.Range("B45").Value = 88.20553994916419
$ws.Range("C45").Value = 34.94953567365734
$ws.Range("B46").Value = 88.40445056698064
$ws.Range("C46").Value = 35.81970242931871
$ws.Range("B47").Value = 98.16248239278272
$ws.Range("C47").Value = 36.60243725342756
